$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 20
$ws.Columns.Item(7).ColumnWidth = 30
$ws.Columns.Item(8).ColumnWidth = 8.43
$ws.Columns.Item(9).ColumnWidth = 18.77734375
